# Corrected some selection scopes
# The data series had too many observations selected; only every other
# quarter (plus the final contiguous block) should have been included.
# Remove the extraneous rows (2,4,6,...,30) so the remaining rows shift
# up and the sheet/shared-strings table are recompacted accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from the bottom up so row numbers of not-yet-deleted rows
# are unaffected by earlier deletions.
$rowsToDelete = @(30, 28, 26, 24, 22, 20, 18, 16, 14, 12, 10, 8, 6, 4, 2)
foreach ($r in $rowsToDelete) {
    $ws.Rows("$r`:$r").Delete()
}
